$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-05-25T17:07:43"
$ws.Range("U4").Value = 35.15
$ws.Range("V4").Value = 39
$ws.Range("W4").Value = 56.8
$ws.Range("X4").Value = 33.41
$ws.Range("Y4").Value = 27.21
$ws.Range("Z4").Value = 7.07
$ws.Range("U6").Value = -0.18
$ws.Range("V6").Value = -0.16
$ws.Range("W6").Value = -0.62
$ws.Range("X6").Value = -0.13
$ws.Range("Y6").Value = 0.05
$ws.Range("Z6").Value = 0.05
$ws.Range("U9").Value = 35.43
$ws.Range("V9").Value = 39.43
$ws.Range("W9").Value = 57.25
$ws.Range("X9").Value = 34.06
$ws.Range("Y9").Value = 27.85
$ws.Range("Z9").Value = 7.23
$ws.Range("U11").Value = 0.11
$ws.Range("V11").Value = 0.28
$ws.Range("W11").Value = -0.17
$ws.Range("X11").Value = 0.51
$ws.Range("Y11").Value = 0.7
$ws.Range("Z11").Value = 0.2
$ws.Range("U14").Value = 35.43
$ws.Range("V14").Value = 39.43
$ws.Range("W14").Value = 57.25
$ws.Range("X14").Value = 34.06
$ws.Range("Y14").Value = 27.85
$ws.Range("Z14").Value = 7.23
$ws.Range("U16").Value = 0.11
$ws.Range("V16").Value = 0.28
$ws.Range("W16").Value = -0.17
$ws.Range("X16").Value = 0.51
$ws.Range("Y16").Value = 0.7
$ws.Range("Z16").Value = 0.2
$ws.Range("U19").Value = 35.4
$ws.Range("V19").Value = 39.35
$ws.Range("W19").Value = 57.19
$ws.Range("X19").Value = 33.68
$ws.Range("Y19").Value = 27.46
$ws.Range("Z19").Value = 7.13
$ws.Range("U21").Value = 0.07000000000000001
$ws.Range("V21").Value = 0.2
$ws.Range("W21").Value = -0.23
$ws.Range("X21").Value = 0.13
$ws.Range("Y21").Value = 0.3
$ws.Range("Z21").Value = 0.11
$ws.Range("U24").Value = 35.4
$ws.Range("V24").Value = 39.35
$ws.Range("W24").Value = 57.19
$ws.Range("X24").Value = 33.68
$ws.Range("Y24").Value = 27.46
$ws.Range("Z24").Value = 7.13
$ws.Range("U26").Value = 0.07000000000000001
$ws.Range("V26").Value = 0.2
$ws.Range("W26").Value = -0.23
$ws.Range("X26").Value = 0.13
$ws.Range("Y26").Value = 0.3
$ws.Range("Z26").Value = 0.11
$ws.Range("V29").Value = 39.87
$ws.Range("W29").Value = 57.77
$ws.Range("Y29").Value = 27.77
$ws.Range("Z29").Value = 7.22
$ws.Range("U31").Value = 0.5
$ws.Range("V31").Value = 0.72
$ws.Range("W31").Value = 0.35
$ws.Range("X31").Value = 0.51
$ws.Range("Z31").Value = 0.19
$ws.Range("U34").Value = 34.7
$ws.Range("V34").Value = 38.73
$ws.Range("W34").Value = 56.6
$ws.Range("X34").Value = 34.69
$ws.Range("Y34").Value = 28.71
$ws.Range("Z34").Value = 7.4
$ws.Range("W35").Value = 0.19
$ws.Range("U36").Value = -0.62
$ws.Range("V36").Value = -0.43
$ws.Range("W36").Value = -1.02
$ws.Range("X36").Value = 1.14
$ws.Range("Y36").Value = 1.55
$ws.Range("Z36").Value = 0.37
$ws.Range("U39").Value = 35.15
$ws.Range("V39").Value = 39
$ws.Range("W39").Value = 56.8
$ws.Range("X39").Value = 33.41
$ws.Range("Y39").Value = 27.21
$ws.Range("Z39").Value = 7.07
$ws.Range("U41").Value = -0.18
$ws.Range("V41").Value = -0.16
$ws.Range("W41").Value = -0.62
$ws.Range("X41").Value = -0.13
$ws.Range("Y41").Value = 0.05
$ws.Range("Z41").Value = 0.05
$ws.Range("U44").Value = 35.83
$ws.Range("V44").Value = 39.67
$ws.Range("W44").Value = 57.89
$ws.Range("X44").Value = 33.92
$ws.Range("Y44").Value = 27.48
$ws.Range("Z44").Value = 7.16
$ws.Range("U46").Value = 0.5
$ws.Range("V46").Value = 0.52
$ws.Range("W46").Value = 0.46
$ws.Range("X46").Value = 0.37
$ws.Range("Y46").Value = 0.33
$ws.Range("Z46").Value = 0.14
$ws.Range("U49").Value = 32.23
$ws.Range("V49").Value = 35.56
$ws.Range("W49").Value = 51.83
$ws.Range("X49").Value = 30.11
$ws.Range("Y49").Value = 24.38
$ws.Range("Z49").Value = 6.62
$ws.Range("U51").Value = -3.09
$ws.Range("V51").Value = -3.59
$ws.Range("W51").Value = -5.6
$ws.Range("X51").Value = -3.43
$ws.Range("Y51").Value = -2.78
$ws.Range("Z51").Value = -0.4
$ws.Range("U54").Value = 30.99
$ws.Range("V54").Value = 34.37
$ws.Range("W54").Value = 48.17
$ws.Range("X54").Value = 28.97
$ws.Range("Y54").Value = 23.41
$ws.Range("Z54").Value = 6.22
$ws.Range("U56").Value = -4.34
$ws.Range("V56").Value = -4.78
$ws.Range("W56").Value = -9.25
$ws.Range("X56").Value = -4.58
$ws.Range("Y56").Value = -3.75
$ws.Range("Z56").Value = -0.8
$ws.Range("U59").Value = 36.8
$ws.Range("V59").Value = 40.66
$ws.Range("W59").Value = 59.32
$ws.Range("X59").Value = 34.73
$ws.Range("Y59").Value = 28.17
$ws.Range("Z59").Value = 7.31
$ws.Range("W61").Value = 1.9
$ws.Range("X61").Value = 1.18
$ws.Range("Y61").Value = 1.01
$ws.Range("Z61").Value = 0.29
$ws.Range("U64").Value = 37.38
$ws.Range("V64").Value = 41.26
$ws.Range("W64").Value = 60.19
$ws.Range("X64").Value = 35.24
$ws.Range("Y64").Value = 28.55
$ws.Range("Z64").Value = 7.4
$ws.Range("U66").Value = 2.06
$ws.Range("V66").Value = 2.1
$ws.Range("W66").Value = 2.77
$ws.Range("Y66").Value = 1.4
$ws.Range("Z66").Value = 0.37
$ws.Range("U69").Value = 37.78
$ws.Range("V69").Value = 41.61
$ws.Range("W69").Value = 60.64
$ws.Range("X69").Value = 35.61
$ws.Range("Y69").Value = 28.92
$ws.Range("Z69").Value = 7.51
$ws.Range("U71").Value = 2.46
$ws.Range("V71").Value = 2.46
$ws.Range("W71").Value = 3.21
$ws.Range("Y71").Value = 1.76
$ws.Range("Z71").Value = 0.48
$ws.Range("U74").Value = 36.09
$ws.Range("V74").Value = 39.83
$ws.Range("W74").Value = 58.12
$ws.Range("X74").Value = 34.02
$ws.Range("Y74").Value = 27.65
$ws.Range("Z74").Value = 7.18
$ws.Range("U76").Value = 0.76
$ws.Range("V76").Value = 0.68
$ws.Range("W76").Value = 0.7
$ws.Range("X76").Value = 0.48
$ws.Range("Y76").Value = 0.5
$ws.Range("Z76").Value = 0.15
$ws.Range("U79").Value = 35.33
$ws.Range("V79").Value = 39.16
$ws.Range("W79").Value = 57.42
$ws.Range("X79").Value = 33.55
$ws.Range("Y79").Value = 27.16
$ws.Range("Z79").Value = 7.03
$ws.Range("U84").Value = 30.99
$ws.Range("V84").Value = 34.38
$ws.Range("W84").Value = 48.17
$ws.Range("X84").Value = 28.97
$ws.Range("Y84").Value = 23.43
$ws.Range("Z84").Value = 6.23
$ws.Range("U86").Value = -4.34
$ws.Range("V86").Value = -4.78
$ws.Range("W86").Value = -9.25
$ws.Range("X86").Value = -4.58
$ws.Range("Y86").Value = -3.73
$ws.Range("Z86").Value = -0.8
$ws.Range("V89").Value = 39.87
$ws.Range("W89").Value = 57.77
$ws.Range("Y89").Value = 27.77
$ws.Range("Z89").Value = 7.22
$ws.Range("U91").Value = 0.5
$ws.Range("V91").Value = 0.72
$ws.Range("W91").Value = 0.35
$ws.Range("X91").Value = 0.51
$ws.Range("Z91").Value = 0.19
